$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.905654311180115
$ws.Range("B1").Value = 4.73459529876709
$ws.Range("C1").Value = 3.319035291671753
$ws.Range("D1").Value = 2.212757110595703
$ws.Range("E1").Value = 1.98274290561676
